# Update daily power records: extend the comforter_cda_table from A1:F59
# to A1:F65 (six new daily rows, 2018-10-10 .. 2018-10-16) and fill in the
# start/end time readings that are known so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert six new rows below the last table row (row 59). Inserting (rather
# than just writing past the end of the table) makes the new D/E/F cells
# pick up the same per-cell style as the row above, matching the existing
# "Duration" / "Second Duration" / "Absolute Value" calculated columns.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows(60 + $i).Insert()
}

# Grow the table definition to cover the newly inserted rows; keeps the
# table ref / autoFilter ref and the sheet dimension in sync with the data.
$lo.Resize($ws.Range("A1:F65"))

# --- Date column (A) ---------------------------------------------------
$ws.Range("A59").Value = 43383
$ws.Range("A60").Value = 43384
$ws.Range("A61").Value = 43385
$ws.Range("A62").Value = 43386
$ws.Range("A63").Value = 43387
$ws.Range("A64").Value = 43388
$ws.Range("A65").Value = 43389

# --- Start / End time columns (B / C) -----------------------------------
# Only row 62 has both readings recorded; row 63 only has a start time so
# far (its end time is still blank, matching an in-progress entry).
$ws.Range("B62").Value = 0.73958333333333337
$ws.Range("C62").Value = 0.99930555555555556
$ws.Range("B63").Value = 0

# --- Calculated columns (D / E / F) -------------------------------------
# Fill each pair of formula ranges together so Excel stores them as
# shared formulas, matching how the table's calculated columns behave.
$ws.Range("D60:D61").Formula = "=(C60-B60)* 1440"
$ws.Range("E60:E61").Formula = "=IF(C60>B60, (C60-B60)*1440, (B60-C60)*1440)"
$ws.Range("F60:F61").Formula = "=ABS((C60-B60)*1440)"

$ws.Range("D62:D65").Formula = "=(C62-B62)* 1440"
$ws.Range("E62:E65").Formula = "=IF(C62>B62, (C62-B62)*1440, (B62-C62)*1440)"
$ws.Range("F62:F65").Formula = "=ABS((C62-B62)*1440)"

# --- View state ----------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("C63").Select() | Out-Null

Write-Host "Updated table range: $($lo.Range.Address())"
